# Apply the "New Orleans shard" edit:
#  1. Insert a new "State" column into hotel_info (between Hotel_Name and City)
#     and populate it with "Louisiana" for the existing data row.
#  2. Re-order the worksheet tabs so "review_info" comes before "hotel_info".

$wb = $excel.ActiveWorkbook

# --- 1. hotel_info: insert State column -------------------------------------
$hotel = $wb.Worksheets.Item("hotel_info")

# Column C currently holds "City"; insert a new blank column there and push
# City (and everything after it) one column to the right.
$hotel.Columns("C").Insert()

$hotel.Range("C1").Value = "State"
$hotel.Range("C2").Value = "Louisiana"

# --- 2. Reorder sheet tabs: review_info first, hotel_info second ------------
$review = $wb.Worksheets.Item("review_info")
$review.Move($hotel)
